# Applies updated probability values to the team-specific matrix sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1923076923076923
$ws.Range("C2").Value = 0.5621301775147929
$ws.Range("J2").Value = 0.008875739644970414
$ws.Range("P2").Value = 0.1538461538461539
$ws.Range("S2").Value = 0.08284023668639054
$ws.Range("B3").Value = 0.01
$ws.Range("C3").Value = 0.04
$ws.Range("J3").Value = 0.005
$ws.Range("P3").Value = 0.8
$ws.Range("S3").Value = 0.145
$ws.Range("J4").Value = 0.04545454545454546
$ws.Range("P4").Value = 0.4318181818181818
$ws.Range("S4").Value = 0.5227272727272727
$ws.Range("B6").Value = 0.06862745098039216
$ws.Range("D6").Value = 0.004901960784313725
$ws.Range("F6").Value = 0.04411764705882353
$ws.Range("J6").Value = 0.2941176470588235
$ws.Range("O6").Value = 0.01470588235294118
$ws.Range("Q6").Value = 0.1127450980392157
$ws.Range("R6").Value = 0.107843137254902
$ws.Range("S6").Value = 0.3529411764705883
$ws.Range("B7").Value = 0.1223404255319149
$ws.Range("D7").Value = 0.01595744680851064
$ws.Range("F7").Value = 0.0425531914893617
$ws.Range("J7").Value = 0.2074468085106383
$ws.Range("Q7").Value = 0.1702127659574468
$ws.Range("R7").Value = 0.0797872340425532
$ws.Range("S7").Value = 0.3617021276595745
$ws.Range("B8").Value = 0.1068702290076336
$ws.Range("D8").Value = 0.01526717557251908
$ws.Range("F8").Value = 0.05089058524173028
$ws.Range("J8").Value = 0.1145038167938931
$ws.Range("O8").Value = 0.01017811704834606
$ws.Range("Q8").Value = 0.1653944020356234
$ws.Range("R8").Value = 0.1272264631043257
$ws.Range("S8").Value = 0.4096692111959288
$ws.Range("B9").Value = 0.08823529411764706
$ws.Range("D9").Value = 0.01176470588235294
$ws.Range("F9").Value = 0.05882352941176471
$ws.Range("J9").Value = 0.08823529411764706
$ws.Range("O9").Value = 0.02352941176470588
$ws.Range("Q9").Value = 0.1764705882352941
$ws.Range("R9").Value = 0.1235294117647059
$ws.Range("S9").Value = 0.4294117647058823
$ws.Range("B10").Value = 0.1329746348962337
$ws.Range("D10").Value = 0.02459646425826287
$ws.Range("E10").Value = 0.0007686395080707148
$ws.Range("F10").Value = 0.07148347425057648
$ws.Range("J10").Value = 0.1145272867025365
$ws.Range("O10").Value = 0.008455034588777863
$ws.Range("Q10").Value = 0.1906225980015373
$ws.Range("R10").Value = 0.09069946195234435
$ws.Range("S10").Value = 0.3658724058416603
$ws.Range("G11").Value = 0.1466666666666667
$ws.Range("K11").Value = 0.2233333333333333
$ws.Range("L11").Value = 0.49
$ws.Range("S11").Value = 0.04
$ws.Range("G12").Value = 0.7337662337662337
$ws.Range("J12").Value = 0.1948051948051948
$ws.Range("K12").Value = 0.006493506493506494
$ws.Range("L12").Value = 0.01948051948051948
$ws.Range("S12").Value = 0.04545454545454546
$ws.Range("G13").Value = 0.6101694915254238
$ws.Range("J13").Value = 0.3220338983050847
$ws.Range("S13").Value = 0.06779661016949153
$ws.Range("F15").Value = 0.03
$ws.Range("H15").Value = 0.125
$ws.Range("I15").Value = 0.05
$ws.Range("J15").Value = 0.41
$ws.Range("K15").Value = 0.09
$ws.Range("M15").Value = 0.02
$ws.Range("O15").Value = 0.025
$ws.Range("S15").Value = 0.25
$ws.Range("F16").Value = 0.01851851851851852
$ws.Range("H16").Value = 0.1712962962962963
$ws.Range("I16").Value = 0.05555555555555555
$ws.Range("J16").Value = 0.412037037037037
$ws.Range("K16").Value = 0.1203703703703704
$ws.Range("M16").Value = 0.02777777777777778
$ws.Range("O16").Value = 0.05092592592592592
$ws.Range("S16").Value = 0.1435185185185185
$ws.Range("F17").Value = 0.01754385964912281
$ws.Range("H17").Value = 0.1854636591478697
$ws.Range("I17").Value = 0.07769423558897243
$ws.Range("J17").Value = 0.4235588972431077
$ws.Range("K17").Value = 0.08771929824561403
$ws.Range("M17").Value = 0.02005012531328321
$ws.Range("N17").Value = 0.002506265664160401
$ws.Range("O17").Value = 0.05263157894736842
$ws.Range("S17").Value = 0.1328320802005012
$ws.Range("F18").Value = 0.02202643171806168
$ws.Range("H18").Value = 0.13215859030837
$ws.Range("I18").Value = 0.08370044052863436
$ws.Range("J18").Value = 0.4273127753303965
$ws.Range("K18").Value = 0.07048458149779736
$ws.Range("M18").Value = 0.013215859030837
$ws.Range("O18").Value = 0.1101321585903084
$ws.Range("S18").Value = 0.1409691629955947
$ws.Range("F19").Value = 0.01610305958132045
$ws.Range("H19").Value = 0.1859903381642512
$ws.Range("I19").Value = 0.07890499194847021
$ws.Range("J19").Value = 0.3864734299516908
$ws.Range("K19").Value = 0.1103059581320451
$ws.Range("M19").Value = 0.0322061191626409
$ws.Range("N19").Value = 0.005636070853462158
$ws.Range("O19").Value = 0.06763285024154589
$ws.Range("S19").Value = 0.1167471819645733
